# "three and four teams tied"
# Recompute change_flag (P) / change_count (Q) running totals for several
# groups now that the set of teams considered tied changed (Belgium ->
# Slovakia swap in the play-off pool, and an adjustment to Ukraine's /
# Slovakia's goal-difference entries), which removes some rows that used
# to register a ranking change.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("Q5").Value = 1
$ws.Range("Q6").Value = 1
$ws.Range("Q7").Value = 1
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 1
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 1
$ws.Range("Q10").Value = 1
$ws.Range("Q11").Value = 1
$ws.Range("Q12").Value = 1
$ws.Range("Q13").Value = 2
$ws.Range("Q14").Value = 2
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 2
$ws.Range("Q16").Value = 2
$ws.Range("Q17").Value = 3
$ws.Range("Q18").Value = 3
$ws.Range("Q19").Value = 4
$ws.Range("Q20").Value = 5
$ws.Range("Q21").Value = 6
$ws.Range("Q22").Value = 7
$ws.Range("Q23").Value = 8
$ws.Range("Q24").Value = 8
$ws.Range("Q25").Value = 9
$ws.Range("Q26").Value = 10
$ws.Range("Q27").Value = 10
$ws.Range("Q28").Value = 10
$ws.Range("Q29").Value = 11
$ws.Range("P36").Value = 0
$ws.Range("Q36").Value = 1
$ws.Range("Q37").Value = 2
$ws.Range("Q38").Value = 2
$ws.Range("Q39").Value = 2
$ws.Range("Q40").Value = 2
$ws.Range("P41").Value = 0
$ws.Range("Q41").Value = 2
$ws.Range("Q42").Value = 3
$ws.Range("P43").Value = 0
$ws.Range("Q43").Value = 3
$ws.Range("P44").Value = 0
$ws.Range("Q44").Value = 3
$ws.Range("Q45").Value = 4
$ws.Range("Q46").Value = 5
$ws.Range("Q47").Value = 5
$ws.Range("Q48").Value = 5
$ws.Range("Q49").Value = 5
$ws.Range("Q50").Value = 5
$ws.Range("Q51").Value = 6
$ws.Range("Q52").Value = 7
$ws.Range("Q53").Value = 8
$ws.Range("Q54").Value = 9
$ws.Range("Q55").Value = 10
$ws.Range("Q56").Value = 10
$ws.Range("Q57").Value = 11
$ws.Range("Q58").Value = 11
$ws.Range("Q59").Value = 11
$ws.Range("Q60").Value = 11
$ws.Range("Q61").Value = 11
$ws.Range("Q62").Value = 11
$ws.Range("Q63").Value = 11
$ws.Range("Q64").Value = 11
$ws.Range("Q65").Value = 11
$ws.Range("P66").Value = 0
$ws.Range("Q66").Value = 11
$ws.Range("Q67").Value = 12
$ws.Range("Q68").Value = 12
$ws.Range("Q69").Value = 12
$ws.Range("Q70").Value = 13
$ws.Range("Q71").Value = 14
$ws.Range("Q72").Value = 15
$ws.Range("Q73").Value = 16
$ws.Range("Q74").Value = 17
$ws.Range("K75").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M75").Value = '[''Austria'', ''Slovakia'', ''Slovenia'', ''Scotland'']'
$ws.Range("K76").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M76").Value = '[''Austria'', ''Slovakia'', ''Slovenia'', ''Scotland'']'
$ws.Range("K77").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M77").Value = '[''Austria'', ''Slovakia'', ''Slovenia'', ''Scotland'']'
$ws.Range("K78").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M78").Value = '[''Austria'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K79").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M79").Value = '[''Austria'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K80").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M80").Value = '[''Austria'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K81").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M81").Value = '[''Austria'', ''Slovakia'', ''Italy'', ''Hungary'']'
$ws.Range("K82").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M82").Value = '[''Austria'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K83").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M83").Value = '[''Austria'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K84").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M84").Value = '[''Netherlands'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K85").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M85").Value = '[''Austria'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K86").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M86").Value = '[''Austria'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K87").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M87").Value = '[''Netherlands'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K88").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M88").Value = '[''Austria'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K89").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M89").Value = '[''Austria'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K90").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M90").Value = '[''Netherlands'', ''Slovakia'', ''Hungary'', ''Slovenia'']'
$ws.Range("K91").Value = '[''Slovakia'', 3, 1, 2]'
$ws.Range("M91").Value = '[''Netherlands'', ''Slovakia'', ''Slovenia'', ''Hungary'']'
$ws.Range("P91").Value = 0
$ws.Range("Q91").Value = 8
$ws.Range("K92").Value = '[''Slovakia'', 4, 1, 2]'
$ws.Range("M92").Value = '[''Slovakia'', ''Netherlands'', ''Slovenia'', ''Hungary'']'
$ws.Range("P92").Value = 0
$ws.Range("Q92").Value = 8
$ws.Range("K93").Value = '[''Ukraine'', 4, 0, 2]'
$ws.Range("Q93").Value = 9
$ws.Range("K94").Value = '[''Slovakia'', 4, 1, 3]'
$ws.Range("M94").Value = '[''Slovakia'', ''Netherlands'', ''Slovenia'', ''Hungary'']'
$ws.Range("Q94").Value = 10
$ws.Range("K95").Value = '[''Slovakia'', 4, 1, 3]'
$ws.Range("M95").Value = '[''Slovakia'', ''Netherlands'', ''Slovenia'', ''Hungary'']'
$ws.Range("Q95").Value = 10
$ws.Range("K96").Value = '[''Slovakia'', 4, 1, 3]'
$ws.Range("M96").Value = '[''Slovakia'', ''Netherlands'', ''Georgia'', ''Slovenia'']'
$ws.Range("Q96").Value = 11
$ws.Range("K97").Value = '[''Slovakia'', 4, 1, 3]'
$ws.Range("M97").Value = '[''Slovakia'', ''Netherlands'', ''Georgia'', ''Slovenia'']'
$ws.Range("Q97").Value = 11
$ws.Range("K98").Value = '[''Slovakia'', 4, 1, 3]'
$ws.Range("M98").Value = '[''Slovakia'', ''Netherlands'', ''Georgia'', ''Slovenia'']'
$ws.Range("P98").Value = 0
$ws.Range("Q98").Value = 11
$ws.Range("K99").Value = '[''Slovakia'', 4, 1, 3]'
$ws.Range("M99").Value = '[''Slovakia'', ''Netherlands'', ''Georgia'', ''Slovenia'']'
$ws.Range("Q99").Value = 11
$ws.Range("K100").Value = '[''Slovakia'', 4, 1, 3]'
$ws.Range("M100").Value = '[''Slovakia'', ''Netherlands'', ''Georgia'', ''Slovenia'']'
$ws.Range("Q100").Value = 11
